# Negate specific cell values to test negative power handling, then
# move the selection to B6 (matches the saved sheetView selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = -80
$ws.Range("A9").Value = -200
$ws.Range("A26").Value = -420
$ws.Range("A28").Value = -60
$ws.Range("B31").Value = -600

$ws.Range("B6").Select()
